# Apply corrected mean_Intake / sem_Intake values (column O/P):
# Intake was previously just milligrams infused (EarnedInfusions * mg/infusion);
# it now also divides by animal weight, so the statistic is recomputed per group/session.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = 1434.2544143502885
$ws.Range("P2").Value = 407.27849891053671
$ws.Range("O3").Value = 1082.0945214997121
$ws.Range("P3").Value = 263.78677607768248
$ws.Range("O4").Value = 938.81770506466887
$ws.Range("P4").Value = 274.01021130109228
$ws.Range("O5").Value = 994.12548142685989
$ws.Range("P5").Value = 218.99220068394027
$ws.Range("O6").Value = 978.83825359412697
$ws.Range("P6").Value = 100.33755111210807
$ws.Range("O7").Value = 1082.2797891025903
$ws.Range("P7").Value = 163.53992172771223
$ws.Range("O8").Value = 1080.6866296899257
$ws.Range("P8").Value = 216.60701676184294
$ws.Range("O9").Value = 1324.8638894545493
$ws.Range("P9").Value = 320.47656209320991
$ws.Range("O10").Value = 1582.4350881425946
$ws.Range("P10").Value = 579.64834031905082
$ws.Range("O11").Value = 1698.6328075881597
$ws.Range("P11").Value = 565.99394444528457
$ws.Range("O12").Value = 1670.8375799678433
$ws.Range("P12").Value = 466.89452912649813
$ws.Range("O13").Value = 1652.8115786520668
$ws.Range("P13").Value = 438.26982011107418
$ws.Range("O14").Value = 1463.6818290485039
$ws.Range("P14").Value = 298.03643026301012
$ws.Range("O15").Value = 1457.1653528181744
$ws.Range("P15").Value = 331.75399626602655
$ws.Range("O16").Value = 1220.0223486657064
$ws.Range("P16").Value = 248.42509807874515
$ws.Range("O28").Value = 78.953200851602801
$ws.Range("P28").Value = 25.361198454273815
$ws.Range("O29").Value = 65.799865891336751
$ws.Range("P29").Value = 28.941686015155629
$ws.Range("O30").Value = 385.79834114126686
$ws.Range("P30").Value = 138.22841683717024
$ws.Range("O31").Value = 204.02210154321506
$ws.Range("P31").Value = 40.748089529932031
$ws.Range("O32").Value = 913.77346131005288
$ws.Range("P32").Value = 304.50916542082626
$ws.Range("O33").Value = 1114.1799491819634
$ws.Range("P33").Value = 239.91631804323137
$ws.Range("O34").Value = 1526.1895540730156
$ws.Range("P34").Value = 337.13181432097474
$ws.Range("O35").Value = 2180.3988014967167
$ws.Range("P35").Value = 894.3400661644464
$ws.Range("O36").Value = 2158.1985583312608
$ws.Range("P36").Value = 658.03523330508551
$ws.Range("O37").Value = 2138.3516309255801
$ws.Range("P37").Value = 489.40000604564307
$ws.Range("O38").Value = 2228.3962112146328
$ws.Range("P38").Value = 506.71802636398144
$ws.Range("O39").Value = 1778.8293491736713
$ws.Range("P39").Value = 561.32813468278528
$ws.Range("O40").Value = 1718.9336679522955
$ws.Range("P40").Value = 397.58916984157867
$ws.Range("O41").Value = 2551.7751498423527
$ws.Range("P41").Value = 559.41546036939576
$ws.Range("O42").Value = 2657.4939003195691
$ws.Range("P42").Value = 678.62124970744344
$ws.Range("O54").Value = 1568.3319974881877
$ws.Range("P54").Value = 459.2497981802731
$ws.Range("O55").Value = 1794.5064494284131
$ws.Range("P55").Value = 833.34057502130474
$ws.Range("O56").Value = 1924.1407394899913
$ws.Range("P56").Value = 937.28132389492339
$ws.Range("O57").Value = 487.37926932367156
$ws.Range("P57").Value = 56.515743989173771
$ws.Range("O58").Value = 1325.2577192843924
$ws.Range("P58").Value = 436.37342576577629
$ws.Range("O59").Value = 990.69535128031532
$ws.Range("P59").Value = 235.43650583171947
$ws.Range("O60").Value = 1075.3275426669525
$ws.Range("P60").Value = 204.62343609831993
$ws.Range("O61").Value = 1102.8672143018716
$ws.Range("P61").Value = 250.84447249522464
$ws.Range("O62").Value = 1054.4674753130462
$ws.Range("P62").Value = 194.0020198432571
$ws.Range("O63").Value = 1095.5945688059478
$ws.Range("P63").Value = 214.80171631264878
$ws.Range("O64").Value = 1605.1071519143266
$ws.Range("P64").Value = 333.20815766953075
$ws.Range("O65").Value = 782.41116158378384
$ws.Range("P65").Value = 147.33079425288483
$ws.Range("O66").Value = 1235.51286702078
$ws.Range("P66").Value = 367.89252947131507
$ws.Range("O67").Value = 902.20216683175363
$ws.Range("P67").Value = 278.44982709496821
$ws.Range("O68").Value = 1378.935793313713
$ws.Range("P68").Value = 353.16146787877869
$ws.Range("O80").Value = 145.15399385742455
$ws.Range("P80").Value = 64.099312183109319
$ws.Range("O81").Value = 191.24440876161913
$ws.Range("P81").Value = 119.42867895551011
$ws.Range("O82").Value = 155.69102944150038
$ws.Range("P82").Value = 40.323196221567471
$ws.Range("O83").Value = 330.32900685237746
$ws.Range("P83").Value = 147.15509549674047
$ws.Range("O84").Value = 553.18523388621361
$ws.Range("P84").Value = 171.67898060644472
$ws.Range("O85").Value = 773.53062596702864
$ws.Range("P85").Value = 224.85153185499897
$ws.Range("O86").Value = 791.8141593991578
$ws.Range("P86").Value = 239.08678493895118
$ws.Range("O87").Value = 833.42070726168208
$ws.Range("P87").Value = 246.72777098168928
$ws.Range("O88").Value = 949.51606851535087
$ws.Range("P88").Value = 309.00998992584141
$ws.Range("O89").Value = 1049.7105752846203
$ws.Range("P89").Value = 338.40303752922205
$ws.Range("O90").Value = 851.55719187298712
$ws.Range("P90").Value = 257.6966510767046
$ws.Range("O91").Value = 520.34415625000008
$ws.Range("P91").Value = 246.91378074268985
$ws.Range("O92").Value = 1248.7215054901369
$ws.Range("P92").Value = 499.71782656036186
$ws.Range("O93").Value = 1543.0769206891837
$ws.Range("P93").Value = 677.73470124161622
$ws.Range("O94").Value = 1082.417938893059
$ws.Range("P94").Value = 352.35404399770698

# Column P (sem_Intake) got slightly narrower in the saved workbook (12.7109375 -> 11.7109375 chars).
# Excel quantizes ColumnWidth to whole pixels, so 10.83 is the closest settable value
# that reproduces the target width exactly (11.666... chars, the nearest achievable pixel width).
$ws.Columns.Item(16).ColumnWidth = 10.83
